$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("model_20_2_0", 0.935668962192915, 0.8269569456985931, 0.9986499951093522, 0.8615366889925381, 0.9750981382103319, 0.4301819460331697, 1.157139700903508, 0.01088443847589244, 0.3827732931438831, 0.1968288658098877, 0.486977646559948, 0.6558825703074977, 1.029131035988114, 0.6838048491895802, 155.6870940589584, 249.5405325738099),
    @("model_20_2_1", 0.9410905465579085, 0.8251714933389779, 0.9939360014962416, 0.7948413273062608, 0.9610307888370628, 0.3939277864203578, 1.169079029053507, 0.04889109594290887, 0.5671485117077592, 0.3080197657305038, 0.7470559421610503, 0.6276366675237814, 1.026675978917174, 0.6543564598472816, 155.8631753393594, 249.7166138542108),
    @("model_20_2_2", 0.9419179992932289, 0.8206058134564078, 0.9872788241544251, 0.6966507878360164, 0.9404647621743746, 0.3883946061692002, 1.199609751451276, 0.1025647068328165, 0.8385902089713351, 0.4705773984257574, 0.9557464079471607, 0.623213130613597, 1.026301283338915, 0.6497446038128468, 155.8914668599407, 249.7449053747922),
    @("model_20_2_3", 0.9397853740844054, 0.8143447370368945, 0.9799339364485995, 0.5815368565901007, 0.9165885569282661, 0.4026554807604392, 1.241477598298168, 0.1617829947814116, 1.156815580220497, 0.659299287500954, 1.120945533750144, 0.6345514011964982, 1.027267000414609, 0.661565568047918, 155.8193479382196, 249.672786453071),
    @("model_20_2_4", 0.9359038276609916, 0.8071423968882662, 0.9726232642485966, 0.4604578792128399, 0.89168669090022, 0.4286113995666853, 1.289639680035772, 0.2207254196048964, 1.491530954018616, 0.8561281868117563, 1.250747509996516, 0.6546841983480931, 1.029024681813891, 0.6825554600548918, 155.6944091981611, 249.5478477130126),
    @("model_20_2_5", 0.9310955446105518, 0.7995049192274981, 0.9657333427461596, 0.3405482891000746, 0.8672039073551009, 0.4607644104650836, 1.340711528320679, 0.2762755344352312, 1.823013628765018, 1.049644581600124, 1.351997857362778, 0.6787962952647013, 1.031202017534844, 0.7076940588561615, 155.5497368136339, 249.4031753284853),
    @("model_20_2_6", 0.9258938997175538, 0.7917685892759123, 0.9594502945087772, 0.2260311841664454, 0.8439736394010907, 0.4955478338159343, 1.392444402329458, 0.3269327227569961, 2.139589110441797, 1.233260864288863, 1.430230782823065, 0.7039515848522072, 1.033557479373183, 0.7339202612000356, 155.4041827868229, 249.2576213016743)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$c]
    }
}

# Rows 6-8 are brand new; copy the formatting used by column A of the existing
# model rows (bold, bordered, centered style) onto the new A6:A8 cells.
$ws.Range("A5").Copy()
$ws.Range("A6:A8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

